$wb = $excel.ActiveWorkbook

# --- Update status text: "Ready for handoff" -> "In Translation" ---
# This string is shared across the Overview sheet (columns E/F, the per-locale
# status columns) and the per-locale detail sheets (column C, "Status").
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Narrow the (now shorter) status columns to match the new content ---
# Target character width is ~13.41. The host's xlsx writer snaps ColumnWidth
# to a pixel grid (multiples of 1/6) at save time, so feed it the value in
# that grid's pre-image that lands closest to 13.4101845877511 once rounded.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
